$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 1428708.9
$ws.Range("I5").Value = 1538567.6
$ws.Range("J5").Value = 545
$ws.Range("K5").Value = 1538567.6
$ws.Range("L5").Value = 545
$ws.Range("M5").Value = -1538452.6
$ws.Range("N5").Value = -775
# Row 10
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = -1086
# Row 12
$ws.Range("H12").Value = 129.8
$ws.Range("I12").Value = 124.75
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 124.75
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = 45.25
$ws.Range("N12").Value = -490
# Row 33
$ws.Range("H33").Value = 597.8
$ws.Range("I33").Value = 330.66666
$ws.Range("K33").Value = 330.66666
$ws.Range("M33").Value = -101.66666
# Row 87
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22496
# Row 90
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -72480
# Row 112
$ws.Range("H112").Value = 2385
$ws.Range("J112").Value = 2385
$ws.Range("L112").Value = 7155
$ws.Range("N112").Value = -9371
# Row 132
$ws.Range("H132").Value = 279090.34
$ws.Range("I132").Value = 1380.8438
$ws.Range("K132").Value = 4142.5314
$ws.Range("M132").Value = -1612.5314
# Row 141
$ws.Range("H141").Value = 3639.2285
$ws.Range("J141").Value = 9819.799999999999
$ws.Range("L141").Value = 29459.4
$ws.Range("N141").Value = -39819.39999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 12
$ws.Range("H12").Value = 2700
$ws.Range("J12").Value = 2700
$ws.Range("L12").Value = 2700
$ws.Range("N12").Value = -3046
# Row 32
$ws.Range("H32").Value = 2294.9583
$ws.Range("I32").Value = 1938.2609
$ws.Range("K32").Value = 1938.2609
$ws.Range("M32").Value = -1651.2609
# Row 47
$ws.Range("H47").Value = 35552
$ws.Range("J47").Value = 35552
$ws.Range("L47").Value = 35552
$ws.Range("N47").Value = -37002
# Row 61
$ws.Range("H61").Value = 5396.6665
$ws.Range("I61").Value = 5356.2
$ws.Range("K61").Value = 5356.2
$ws.Range("M61").Value = -5144.2
# Row 74
$ws.Range("H74").Value = 1883.1333
$ws.Range("I74").Value = 2463.2222
$ws.Range("K74").Value = 2463.2222
$ws.Range("M74").Value = -1589.2222
# Row 77
$ws.Range("H77").Value = 1883.1333
$ws.Range("I77").Value = 2463.2222
$ws.Range("K77").Value = 12316.111
$ws.Range("M77").Value = -7948.111000000001
# Row 110
$ws.Range("H110").Value = 1424.5333
$ws.Range("I110").Value = 1368.1111
$ws.Range("K110").Value = 1368.1111
$ws.Range("M110").Value = 676.8888999999999
# Row 136
$ws.Range("H136").Value = 5396.6665
$ws.Range("I136").Value = 5356.2
$ws.Range("K136").Value = 16068.6
$ws.Range("M136").Value = -13518.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 30009
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 30009
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 30009
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = -30503
# Row 64
$ws.Range("H64").Value = 1698
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
# Row 67
$ws.Range("H67").Value = 1698
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
# Row 81
$ws.Range("H81").Value = 44894.668
$ws.Range("J81").Value = 44894.668
$ws.Range("L81").Value = 44894.668
$ws.Range("N81").Value = -47016.668
# Row 84
$ws.Range("H84").Value = 44894.668
$ws.Range("J84").Value = 44894.668
$ws.Range("L84").Value = 134684.004
$ws.Range("N84").Value = -145292.004
# Row 107
$ws.Range("H107").Value = 1204.5714
$ws.Range("I107").Value = 1204.5714
$ws.Range("K107").Value = 1204.5714
$ws.Range("M107").Value = 715.4286
# Row 134
$ws.Range("H134").Value = 1760.4546
$ws.Range("I134").Value = 1760.4546
$ws.Range("K134").Value = 5281.3638
$ws.Range("M134").Value = -2746.3638

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3284.889
$ws.Range("I31").Value = 2431.5
$ws.Range("K31").Value = 2431.5
$ws.Range("M31").Value = -2136.5
# Row 34
$ws.Range("H34").Value = 3284.889
$ws.Range("I34").Value = 2431.5
$ws.Range("K34").Value = 2431.5
$ws.Range("M34").Value = -2229.5
# Row 59
$ws.Range("H59").Value = 44957.383
$ws.Range("J59").Value = 44954
$ws.Range("L59").Value = 44954
$ws.Range("N59").Value = -47244
# Row 99
$ws.Range("H99").Value = 2541
$ws.Range("I99").Value = 2475.4285
$ws.Range("K99").Value = 2475.4285
$ws.Range("M99").Value = -977.4285
# Row 122
$ws.Range("H122").Value = 3061.8667
$ws.Range("J122").Value = 4267.7144
$ws.Range("L122").Value = 12803.1432
$ws.Range("N122").Value = -17703.1432
# Row 126
$ws.Range("H126").Value = 2541
$ws.Range("I126").Value = 2475.4285
$ws.Range("K126").Value = 7426.2855
$ws.Range("M126").Value = -4956.2855

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 362.33334
$ws.Range("I14").Value = 362.33334
$ws.Range("K14").Value = 1087.00002
$ws.Range("M14").Value = -914.0000199999999
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null
# Row 98
$ws.Range("H98").Value = 603.5
$ws.Range("J98").Value = 603.5
$ws.Range("L98").Value = 1810.5
$ws.Range("N98").Value = -4806.5
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").Value = $null
# Row 131
$ws.Range("H131").Value = 31249.684
$ws.Range("J131").Value = 2181.5454
$ws.Range("L131").Value = 6544.6362
$ws.Range("N131").Value = -16624.6362
# Row 140
$ws.Range("H140").Value = 1759.1765
$ws.Range("I140").Value = 1171.9286
$ws.Range("K140").Value = 3515.7858
$ws.Range("M140").Value = 1664.2142

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 10114.625
$ws.Range("I70").Value = 9168
$ws.Range("J70").Value = 10249.857
$ws.Range("K70").Value = 9168
$ws.Range("L70").Value = 10249.857
$ws.Range("M70").Value = -8898
$ws.Range("N70").Value = -10789.857
# Row 73
$ws.Range("H73").Value = 10114.625
$ws.Range("I73").Value = 9168
$ws.Range("J73").Value = 10249.857
$ws.Range("K73").Value = 9168
$ws.Range("L73").Value = 10249.857
$ws.Range("M73").Value = -8232
$ws.Range("N73").Value = -12121.857
# Row 113
$ws.Range("H113").Value = 2252.75
$ws.Range("I113").Value = 1005.5
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1005.5
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 1164.5
$ws.Range("N113").Value = -7840
# Row 132
$ws.Range("H132").Value = 2346.889
$ws.Range("I132").Value = 1589
$ws.Range("K132").Value = 4767
$ws.Range("M132").Value = -2237

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2823.1333
$ws.Range("I7").Value = 1807.5714
$ws.Range("K7").Value = 1807.5714
$ws.Range("M7").Value = -1695.5714
# Row 40
$ws.Range("H40").Value = 2985.2307
$ws.Range("I40").Value = 3025.75
$ws.Range("J40").Value = 2499
$ws.Range("K40").Value = 3025.75
$ws.Range("L40").Value = 2499
$ws.Range("M40").Value = -2889.75
$ws.Range("N40").Value = -2771
# Row 126
$ws.Range("H126").Value = 2823.1333
$ws.Range("I126").Value = 1807.5714
$ws.Range("K126").Value = 5422.7142
$ws.Range("M126").Value = -2952.7142
# Row 132
$ws.Range("H132").Value = 2225.5173
$ws.Range("I132").Value = 2061.5454
$ws.Range("K132").Value = 6184.6362
$ws.Range("M132").Value = -3654.6362

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 4801.1025
$ws.Range("I122").Value = 5529.9165
$ws.Range("J122").Value = 3635
$ws.Range("K122").Value = 16589.7495
$ws.Range("L122").Value = 10905
$ws.Range("M122").Value = -14139.7495
$ws.Range("N122").Value = -15805
# Row 132
$ws.Range("H132").Value = 1918
$ws.Range("I132").Value = 1897.75
$ws.Range("K132").Value = 5693.25
$ws.Range("M132").Value = -3163.25
# Row 136
$ws.Range("H136").Value = 1331.1052
$ws.Range("J136").Value = 1844.3334
$ws.Range("L136").Value = 5533.0002
$ws.Range("N136").Value = -10633.0002
